# Update the timestamp column (Z) to reflect the re-run of the
# pcsmote logging pass recorded on 2025-11-13 (commit: "dataset Us Crime agregado").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("Z2:Z11").Value = "2025-11-13T06:53:05.034869"
$ws.Range("Z12:Z21").Value = "2025-11-13T06:53:05.035867"
$ws.Range("Z22:Z32").Value = "2025-11-13T06:53:05.036866"
$ws.Range("Z33:Z42").Value = "2025-11-13T06:53:05.038126"
$ws.Range("Z43:Z45").Value = "2025-11-13T06:53:05.039136"
$ws.Range("Z46:Z74").Value = "2025-11-13T06:53:05.219983"
$ws.Range("Z75:Z76").Value = "2025-11-13T06:53:05.460467"
$ws.Range("Z77").Value = "2025-11-13T06:53:05.461469"
$ws.Range("Z78:Z79").Value = "2025-11-13T06:53:05.461859"
$ws.Range("Z80:Z81").Value = "2025-11-13T06:53:05.462588"
$ws.Range("Z82:Z85").Value = "2025-11-13T06:53:05.463134"
$ws.Range("Z86:Z95").Value = "2025-11-13T06:53:05.464117"
$ws.Range("Z96").Value = "2025-11-13T06:53:05.465114"
$ws.Range("Z97:Z99").Value = "2025-11-13T06:53:05.466638"
$ws.Range("Z100:Z102").Value = "2025-11-13T06:53:05.467648"
